$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Control 47)
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = [double]"1.942737351341484E-07"
$ws.Range("E3").Value = [double]"1.942737351341484E-07"

# Row 4 (Control 15)
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = [double]"0.999300195218732"
$ws.Range("E4").Value = [double]"0.999300195218732"

# Row 5 (Control 7)
$ws.Range("D5").Value = [double]"0.1444795114167273"
$ws.Range("E5").Value = [double]"0.1444795114167273"

# Row 6 (MDD 43)
$ws.Range("D6").Value = [double]"1.964375447706767E-31"
$ws.Range("E6").Value = [double]"1.964375447706767E-31"

# Row 7 (MDD 3)
$ws.Range("D7").Value = [double]"0.9999999625507473"
$ws.Range("E7").Value = [double]"3.744925269888455E-08"

# Row 8 (MDD 19)
$ws.Range("D8").Value = [double]"0.9999999999999474"
$ws.Range("E8").Value = [double]"5.262457136723242E-14"

# Row 9 (MDD 7)
$ws.Range("D9").Value = [double]"0.9999999999999278"
$ws.Range("E9").Value = [double]"7.216449660063518E-14"

# Row 11 (MDD 1)
$ws.Range("F11").Value = [double]"13.18059539794922"
